$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new header columns for the four door directions
$ws.Range("D1").Value = "North"
$ws.Range("E1").Value = "East"
$ws.Range("F1").Value = "South"
$ws.Range("G1").Value = "West"

# Door data (1 = door present, 0 = no door) for each room row, columns D:G = North/East/South/West
$doorData = @(
    @(1,1,0,0),
    @(1,1,1,0),
    @(1,1,1,0),
    @(1,1,1,0),
    @(1,1,1,0),
    @(1,1,0,1),
    @(0,1,1,0),
    @(1,1,1,0),
    @(0,0,0,1),
    @(1,0,1,1),
    @(1,0,1,1),
    @(1,1,1,1),
    @(0,0,0,1),
    @(0,0,1,1),
    @(1,0,1,0),
    @(0,0,0,1)
)

for ($i = 0; $i -lt $doorData.Count; $i++) {
    $row = $i + 2
    $vals = $doorData[$i]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
    $ws.Range("F$row").Value = $vals[2]
    $ws.Range("G$row").Value = $vals[3]
}

# Update selection to match the committed workbook state
$ws.Range("J22").Select()
